$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.603.18'
$ws.Range("E2").Value = '  +1.03%  '

$ws.Range("D3").Value = '''3.093.19'
$ws.Range("E3").Value = '  +0.33%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '''541.16'
$ws.Range("E5").Value = '  -2.20%  '

$ws.Range("D6").Value = '''136.57'
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '''3.090.01'
$ws.Range("E8").Value = '  +0.52%  '

$ws.Range("E9").Value = '  +1.23%  '

$ws.Range("D10").Value = '''0.156'
$ws.Range("E10").Value = '  -1.10%  '

$ws.Range("D11").Value = '''6.12'
$ws.Range("E11").Value = '  -7.09%  '

$ws.Range("D12").Value = '''0.460'
$ws.Range("E12").Value = '  +2.03%  '

$ws.Range("D13").Value = '''0.0000227'
$ws.Range("E13").Value = '  +5.28%  '

$ws.Range("D14").Value = '''34.84'
$ws.Range("E14").Value = '  +0.12%  '

$ws.Range("D15").Value = '''3.596.24'
$ws.Range("E15").Value = '  +0.43%  '

$ws.Range("D16").Value = '''63.639.41'
$ws.Range("E16").Value = '  +0.92%  '

$ws.Range("E17").Value = '  +0.46%  '

$ws.Range("D18").Value = '''3.094.34'
$ws.Range("E18").Value = '  +0.44%  '

$ws.Range("D19").Value = '''6.72'
$ws.Range("E19").Value = '  +1.42%  '

$ws.Range("D20").Value = '''489.46'
$ws.Range("E20").Value = '  -2.25%  '

$ws.Range("D21").Value = '''13.52'
$ws.Range("E21").Value = '  +0.65%  '

$ws.Range("D22").Value = '''0.704'
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("D23").Value = '''7.23'
$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("D24").Value = '''79.88'
$ws.Range("E24").Value = '  +3.52%  '

$ws.Range("D25").Value = '''12.29'
$ws.Range("E25").Value = '  +0.82%  '

$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").Value = '''2.74'
$ws.Range("E27").Value = '  -0.43%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''8.36'
$ws.Range("E28").Value = '  +3.29%  '

$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.32%  '

$ws.Range("D30").Value = '''26.34'
$ws.Range("E30").Value = '  +0.77%  '

$ws.Range("E31").Value = '  -1.71%  '

$ws.Range("E32").Value = '  +1.28%  '

$ws.Range("D33").Value = '''2.42'
$ws.Range("E33").Value = '  -3.00%  '

$ws.Range("D34").Value = '''57.33'
$ws.Range("E34").Value = '  -1.99%  '

$ws.Range("D35").Value = '''5.41'
$ws.Range("E35").Value = '  +5.34%  '

$ws.Range("D36").Value = '''6.09'
$ws.Range("E36").Value = '  +4.15%  '

$ws.Range("D37").Value = '''492.73'
$ws.Range("E37").Value = '  -6.05%  '

$ws.Range("D38").Value = '''3.198.39'
$ws.Range("E38").Value = '  +5.23%  '

$ws.Range("D39").Value = '''0.0401'
$ws.Range("E39").Value = '  -2.15%  '

$ws.Range("D40").Value = '''0.0805'
$ws.Range("E40").Value = '  +2.59%  '

$ws.Range("E41").Value = '  -2.11%  '

$ws.Range("D42").Value = '''2.71'
$ws.Range("E42").Value = '  +4.23%  '

$ws.Range("D43").Value = '''8.18'
$ws.Range("E43").Value = '  +2.03%  '

$ws.Range("D44").Value = '''0.257'
$ws.Range("E44").Value = '  +2.25%  '


$ws.Range("D46").Value = '''0.0₃0545'
$ws.Range("E46").Value = '  +9.36%  '

$ws.Range("D47").Value = '''2.06'
$ws.Range("E47").Value = '  +0.73%  '

$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '''121.46'
$ws.Range("E48").Value = '  +0.38%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '''24.83'
$ws.Range("E49").Value = '  +5.41%  '

$ws.Range("E50").Value = '  +3.38%  '

$ws.Range("D51").Value = '''2.35'
$ws.Range("E51").Value = '  -0.01%  '
